$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.421.44'
$ws.Range('E2').Value = '  -2.69%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.554.21'
$ws.Range('E3').Value = '  -4.48%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '546.57'
$ws.Range('E5').Value = '  +0.02%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.48'
$ws.Range('E6').Value = '  -3.41%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.585'
$ws.Range('E8').Value = '  -0.14%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.102'
$ws.Range('E9').Value = '  -2.49%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  -0.76%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.42'
$ws.Range('E11').Value = '  +4.45%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.362'
$ws.Range('E12').Value = '  -1.33%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.010.27'
$ws.Range('E13').Value = '  -4.40%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.10'
$ws.Range('E14').Value = '  -4.00%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '61.366.79'
$ws.Range('E15').Value = '  -2.59%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000142'
$ws.Range('E16').Value = '  -1.70%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.570.07'
$ws.Range('E17').Value = '  -4.09%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.44'
$ws.Range('E18').Value = '  -4.11%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.47'
$ws.Range('E19').Value = '  -1.54%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.09'
$ws.Range('E20').Value = '  -2.00%  '

# Row 21
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.997'
$ws.Range('E21').Value = '  +0.25%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.93'
$ws.Range('E22').Value = '  -5.48%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.484'
$ws.Range('E23').Value = '  -3.73%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.89'
$ws.Range('E24').Value = '  -1.09%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.165'
$ws.Range('E25').Value = '  -1.57%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.29%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.03'
$ws.Range('E27').Value = '  -0.90%  '

# Row 28
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.10'
$ws.Range('E28').Value = '  +1.92%  '

# Row 29
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0814'
$ws.Range('E29').Value = '  -3.70%  '

# Row 30
$ws.Range('E30').Value = '  -0.76%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.87'
$ws.Range('E31').Value = '  -2.64%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.21'
$ws.Range('E32').Value = '  -2.78%  '

# Row 33
$ws.Range('E33').Value = '  +0.06%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.78'
$ws.Range('E34').Value = '  +0.12%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.89'
$ws.Range('E35').Value = '  -2.87%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.38'
$ws.Range('E36').Value = '  -2.70%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.76'
$ws.Range('E37').Value = '  -0.09%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.95'
$ws.Range('E38').Value = '  -2.44%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '321.24'
$ws.Range('E39').Value = '  -4.78%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.880'
$ws.Range('E40').Value = '  -5.78%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.86'
$ws.Range('E41').Value = '  -1.31%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.13'
$ws.Range('E42').Value = '  -2.37%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.13%  '

# Row 44
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.37'
$ws.Range('E44').Value = '  -1.23%  '

# Row 45
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.93'
$ws.Range('E45').Value = '  -1.16%  '

# Row 46
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.599'
$ws.Range('E46').Value = '  -2.57%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0956'
$ws.Range('E47').Value = '  -1.24%  '

# Row 48
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0536'
$ws.Range('E48').Value = '  -3.84%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.17'
$ws.Range('E49').Value = '  -5.01%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0235'
$ws.Range('E50').Value = '  -1.38%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.031.86'
$ws.Range('E51').Value = '  -2.42%  '
